$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: set B (date) column to Text format first so the
# ISO date strings are preserved literally instead of being parsed into serial dates ---
$wsForecast.Range("B2:B17").NumberFormat = "@"

$wsForecast.Range("B2").Value = "2025-02-02"
$wsForecast.Range("D2").Value = 67
$wsForecast.Range("E2").Value = 76
$wsForecast.Range("F2").Value = 90
$wsForecast.Range("G2").Value = 105
$wsForecast.Range("H2").Value = 127

$wsForecast.Range("B3").Value = "2025-02-09"
$wsForecast.Range("D3").Value = 68
$wsForecast.Range("E3").Value = 78
$wsForecast.Range("F3").Value = 94
$wsForecast.Range("G3").Value = 110
$wsForecast.Range("H3").Value = 137

$wsForecast.Range("B4").Value = "2025-02-16"
$wsForecast.Range("D4").Value = 68
$wsForecast.Range("E4").Value = 79
$wsForecast.Range("F4").Value = 96
$wsForecast.Range("G4").Value = 114
$wsForecast.Range("H4").Value = 142

$wsForecast.Range("B5").Value = "2025-02-23"
$wsForecast.Range("D5").Value = 68
$wsForecast.Range("E5").Value = 79
$wsForecast.Range("F5").Value = 96
$wsForecast.Range("G5").Value = 114
$wsForecast.Range("H5").Value = 142

$wsForecast.Range("B6").Value = "2025-03-02"
$wsForecast.Range("D6").Value = 68
$wsForecast.Range("E6").Value = 81
$wsForecast.Range("F6").Value = 99
$wsForecast.Range("G6").Value = 119
$wsForecast.Range("H6").Value = 151

$wsForecast.Range("B7").Value = "2025-03-09"
$wsForecast.Range("D7").Value = 68
$wsForecast.Range("E7").Value = 79
$wsForecast.Range("F7").Value = 96
$wsForecast.Range("G7").Value = 116
$wsForecast.Range("H7").Value = 148

$wsForecast.Range("B8").Value = "2025-03-16"
$wsForecast.Range("D8").Value = 67
$wsForecast.Range("E8").Value = 75
$wsForecast.Range("F8").Value = 91
$wsForecast.Range("G8").Value = 112
$wsForecast.Range("H8").Value = 146

$wsForecast.Range("B9").Value = "2025-03-23"
$wsForecast.Range("D9").Value = 67
$wsForecast.Range("E9").Value = 76
$wsForecast.Range("F9").Value = 93
$wsForecast.Range("G9").Value = 115
$wsForecast.Range("H9").Value = 150

$wsForecast.Range("B10").Value = "2025-03-30"
$wsForecast.Range("D10").Value = 67
$wsForecast.Range("E10").Value = 76
$wsForecast.Range("F10").Value = 92
$wsForecast.Range("G10").Value = 111
$wsForecast.Range("H10").Value = 142

$wsForecast.Range("B11").Value = "2025-04-06"
$wsForecast.Range("D11").Value = 66
$wsForecast.Range("E11").Value = 71
$wsForecast.Range("F11").Value = 87
$wsForecast.Range("G11").Value = 107
$wsForecast.Range("H11").Value = 140

$wsForecast.Range("B12").Value = "2025-04-13"
$wsForecast.Range("D12").Value = 65
$wsForecast.Range("E12").Value = 72
$wsForecast.Range("F12").Value = 88
$wsForecast.Range("G12").Value = 110
$wsForecast.Range("H12").Value = 144

$wsForecast.Range("B13").Value = "2025-04-20"
$wsForecast.Range("D13").Value = 63
$wsForecast.Range("E13").Value = 71
$wsForecast.Range("F13").Value = 87
$wsForecast.Range("G13").Value = 108
$wsForecast.Range("H13").Value = 142

$wsForecast.Range("B14").Value = "2025-04-27"
$wsForecast.Range("D14").Value = 62
$wsForecast.Range("E14").Value = 71
$wsForecast.Range("F14").Value = 87
$wsForecast.Range("G14").Value = 107
$wsForecast.Range("H14").Value = 139

$wsForecast.Range("B15").Value = "2025-05-04"
$wsForecast.Range("D15").Value = 62
$wsForecast.Range("E15").Value = 67
$wsForecast.Range("F15").Value = 82
$wsForecast.Range("G15").Value = 102
$wsForecast.Range("H15").Value = 134

$wsForecast.Range("B16").Value = "2025-05-11"
$wsForecast.Range("D16").Value = 61
$wsForecast.Range("E16").Value = 68
$wsForecast.Range("F16").Value = 83
$wsForecast.Range("G16").Value = 102
$wsForecast.Range("H16").Value = 133

$wsForecast.Range("B17").Value = "2025-05-18"
$wsForecast.Range("D17").Value = 60
$wsForecast.Range("E17").Value = 67
$wsForecast.Range("F17").Value = 82
$wsForecast.Range("G17").Value = 101
$wsForecast.Range("H17").Value = 133

# --- Summary sheet updates ---
# Cells whose new value looks numeric or date-like need Text format first so they
# stay as text, matching the workbooks existing inline-string convention.
foreach ($addr in @("B2","B5","B6","B8","B9","B10","B11","B12","B13","B14","B15")) {
    $wsSummary.Range($addr).NumberFormat = "@"
}

$wsSummary.Range("B2").Value = "2022-12-25 to 2025-01-26"
$wsSummary.Range("B5").Value = "149"
$wsSummary.Range("B6").Value = "119"
$wsSummary.Range("B8").Value = "16227 units"
$wsSummary.Range("B9").Value = "1047"
$wsSummary.Range("B10").Value = "541"
$wsSummary.Range("B11").Value = "271"
$wsSummary.Range("B12").Value = "68"
$wsSummary.Range("B13").Value = "2025-03-02"
$wsSummary.Range("B14").Value = "60"
$wsSummary.Range("B15").Value = "2025-05-18"
